$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.368.90'
$ws.Range("E2").Value = '  +1.12%  '

$ws.Range("D3").Value = '3.738.49'
$ws.Range("E3").Value = '  -0.73%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '613.28'
$ws.Range("E5").Value = '  +4.91%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '193.16'
$ws.Range("E6").Value = '  +8.85%  '

$ws.Range("E7").Value = '  +0.20%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("E9").Value = '  +1.29%  '

$ws.Range("E10").Value = '  -3.33%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '60.44'
$ws.Range("E11").Value = '  +12.32%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000293'
$ws.Range("E12").Value = '  -3.50%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.67'
$ws.Range("E13").Value = '  -1.55%  '

$ws.Range("D14").Value = '4.330.89'
$ws.Range("E14").Value = '  -0.65%  '

$ws.Range("D15").Value = '3.731.40'
$ws.Range("E15").Value = '  -0.79%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.58'
$ws.Range("E16").Value = '  -0.54%  '

$ws.Range("E17").Value = '  +0.54%  '

$ws.Range("E18").Value = '  -0.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.99'
$ws.Range("E19").Value = '  -1.37%  '

$ws.Range("D20").Value = '69.187.83'
$ws.Range("E20").Value = '  +0.85%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '414.62'
$ws.Range("E21").Value = '  +0.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.60'
$ws.Range("E22").Value = '  -0.07%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '90.09'
$ws.Range("E23").Value = '  +0.41%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.09'
$ws.Range("E24").Value = '  -0.98%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.98'
$ws.Range("E25").Value = '  -0.29%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.12'
$ws.Range("E26").Value = '  +2.42%  '

$ws.Range("E27").Value = '  -1.52%  '

$ws.Range("E28").Value = '  +1.13%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.78'
$ws.Range("E29").Value = '  +0.84%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.30'
$ws.Range("E30").Value = '  -0.26%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.86'
$ws.Range("E31").Value = '  -1.91%  '

$ws.Range("E32").Value = '  -0.40%  '

$ws.Range("E33").Value = '  +4.30%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '45.89'
$ws.Range("E34").Value = '  +3.43%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '635.66'
$ws.Range("E35").Value = '  +2.92%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '66.96'
$ws.Range("E36").Value = '  +2.28%  '

$ws.Range("D37").Value = '0.0₃0846'
$ws.Range("E37").Value = '  -9.91%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.419'
$ws.Range("E38").Value = '  +2.59%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.15%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.142'
$ws.Range("E41").Value = '  +3.07%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.11'
$ws.Range("E42").Value = '  +0.20%  '

$ws.Range("E43").Value = '  +0.42%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.66'
$ws.Range("E44").Value = '  +0.85%  '

$ws.Range("E45").Value = '  +3.03%  '

$ws.Range("D46").Value = '2.899.24'
$ws.Range("E46").Value = '  +3.65%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.24'
$ws.Range("E47").Value = '  -2.63%  '

$ws.Range("E48").Value = '  +0.15%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '144.21'
$ws.Range("E49").Value = '  +0.58%  '

$ws.Range("E50").Value = '  -1.49%  '

$ws.Range("E51").Value = '  +0.13%  '
